# Update extrapolation calibration outputs (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN)
# Removing less than USD 5 price from extrapolation calibration because it is just a noise
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 118275.4751848024
$ws.Range("E2").Value = -0.02091100561556052
$ws.Range("F2").Value = 0.174796251227516
$ws.Range("G2").Value = -1.649719457582301
$ws.Range("H2").Value = 16.10310679048877
$ws.Range("D3").Value = 119030.1851993516
$ws.Range("E3").Value = -0.02076647844682017
$ws.Range("F3").Value = 0.2050012169011446
$ws.Range("G3").Value = -1.425902590746017
$ws.Range("H3").Value = 16.98311454636362
$ws.Range("D5").Value = 120689.4060706475
$ws.Range("E5").Value = -0.02794997294275259
$ws.Range("F5").Value = 0.2332690681446356
$ws.Range("G5").Value = -0.8338939102106074
$ws.Range("H5").Value = 8.696484875731587
$ws.Range("D6").Value = 121191.8979872002
$ws.Range("E6").Value = -0.03674033151067763
$ws.Range("F6").Value = 0.2572086377334402
$ws.Range("G6").Value = -1.132427611423583
$ws.Range("H6").Value = 9.867556848219831
$ws.Range("D8").Value = 122862.0821874668
$ws.Range("E8").Value = -0.05060856669362174
$ws.Range("F8").Value = 0.2199553743286014
$ws.Range("G8").Value = -0.9051817561684483
$ws.Range("H8").Value = 7.08722495857673
$ws.Range("D9").Value = 124423.6567373964
$ws.Range("E9").Value = -0.09932857608253293
$ws.Range("F9").Value = 0.4512049289244738
$ws.Range("G9").Value = -2.521182240538921
$ws.Range("H9").Value = 14.36116839186669
$ws.Range("D10").Value = 125742.0904017959
$ws.Range("E10").Value = -0.1191110171165711
$ws.Range("F10").Value = 0.4420726491789731
$ws.Range("G10").Value = -1.932609514947561
$ws.Range("H10").Value = 9.879244624286544
$ws.Range("D11").Value = 127827.2940780794
$ws.Range("E11").Value = -0.1966893939678039
$ws.Range("F11").Value = 0.7832383037921936
$ws.Range("G11").Value = -2.650179520439601
$ws.Range("H11").Value = 13.05285757099681
$ws.Range("D13").Value = 118324.6456458961
$ws.Range("E13").Value = -0.0007790735472777527
$ws.Range("F13").Value = 0.1478625512276133
$ws.Range("G13").Value = -0.8083085152215409
$ws.Range("H13").Value = 11.02642042127016
$ws.Range("D18").Value = 118309.8368363478
$ws.Range("E18").Value = -0.003974876193995315
$ws.Range("F18").Value = 0.1395575970205
$ws.Range("G18").Value = -0.9025127284164005
$ws.Range("H18").Value = 10.0909505284078
$ws.Range("D19").Value = 118343.5400836522
$ws.Range("E19").Value = -0.007684672861838372
$ws.Range("F19").Value = 0.1379176838482369
$ws.Range("G19").Value = -0.7507533802123543
$ws.Range("H19").Value = 9.001744558415018
$ws.Range("D20").Value = 118414.1613449437
$ws.Range("E20").Value = -0.01246987575532709
$ws.Range("F20").Value = 0.1564218073129512
$ws.Range("G20").Value = -0.4390436630479489
$ws.Range("H20").Value = 6.469809024023105
